# Auto-generated edit script
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3061.5264
$ws.Range("I40").Value = 6157.5
$ws.Range("J40").Value = 2235.9333
$ws.Range("K40").Value = 6157.5
$ws.Range("L40").Value = 2235.9333
$ws.Range("M40").Value = -5982.5
$ws.Range("N40").Value = -2585.9333

$ws.Range("H86").Value = 280966.5
$ws.Range("I86").Value = 447025.6
$ws.Range("J86").Value = 4201.3335
$ws.Range("K86").Value = 447025.6
$ws.Range("L86").Value = 4201.3335
$ws.Range("M86").Value = -445902.6
$ws.Range("N86").Value = -6447.3335

$ws.Range("H89").Value = 280966.5
$ws.Range("I89").Value = 447025.6
$ws.Range("J89").Value = 4201.3335
$ws.Range("K89").Value = 2235128
$ws.Range("L89").Value = 21006.6675
$ws.Range("M89").Value = -2229512
$ws.Range("N89").Value = -32238.6675

$ws.Range("H138").Value = 6165.3486
$ws.Range("I138").Value = 3362.25
$ws.Range("J138").Value = 6619.9053
$ws.Range("K138").Value = 10086.75
$ws.Range("L138").Value = 19859.7159
$ws.Range("M138").Value = -4946.75
$ws.Range("N138").Value = -30139.7159

$ws.Range("H141").Value = 3384.4546
$ws.Range("I141").Value = 2182
$ws.Range("K141").Value = 6546
$ws.Range("M141").Value = -1366

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2772.7334
$ws.Range("I61").Value = 2720.0908
$ws.Range("J61").Value = 2917.5
$ws.Range("K61").Value = 2720.0908
$ws.Range("L61").Value = 2917.5
$ws.Range("M61").Value = -2508.0908
$ws.Range("N61").Value = -3341.5

$ws.Range("H97").Value = 3345.6667
$ws.Range("I97").Value = 3124.4443
$ws.Range("J97").Value = 4673
$ws.Range("K97").Value = 3124.4443
$ws.Range("L97").Value = 4673
$ws.Range("M97").Value = -2628.4443
$ws.Range("N97").Value = -5665

$ws.Range("H135").Value = 2760435.2
$ws.Range("J135").Value = 2760435.2
$ws.Range("L135").Value = 2760435.2
$ws.Range("N135").Value = -2770575.2

$ws.Range("H136").Value = 2772.7334
$ws.Range("I136").Value = 2720.0908
$ws.Range("J136").Value = 2917.5
$ws.Range("K136").Value = 8160.2724
$ws.Range("L136").Value = 8752.5
$ws.Range("M136").Value = -5610.2724
$ws.Range("N136").Value = -13852.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 1825.6
$ws.Range("I25").Value = 676
$ws.Range("J25").Value = 3550
$ws.Range("K25").Value = 676
$ws.Range("L25").Value = 3550
$ws.Range("M25").Value = -441
$ws.Range("N25").Value = -4020

$ws.Range("H105").Value = 2021.5
$ws.Range("I105").Value = 991.8
$ws.Range("J105").Value = 2757
$ws.Range("K105").Value = 991.8
$ws.Range("L105").Value = 2757
$ws.Range("M105").Value = 755.2
$ws.Range("N105").Value = -6251

$ws.Range("H134").Value = 436565.8
$ws.Range("I134").Value = 556745.25
$ws.Range("J134").Value = 3920
$ws.Range("K134").Value = 1670235.75
$ws.Range("L134").Value = 11760
$ws.Range("M134").Value = -1667700.75
$ws.Range("N134").Value = -16830

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2981.6191
$ws.Range("I16").Value = 3293.4614
$ws.Range("J16").Value = 2474.875
$ws.Range("K16").Value = 3293.4614
$ws.Range("L16").Value = 2474.875
$ws.Range("M16").Value = -3006.4614
$ws.Range("N16").Value = -3048.875

$ws.Range("H31").Value = 2598.0408
$ws.Range("I31").Value = 1392.7407
$ws.Range("J31").Value = 4077.2727
$ws.Range("K31").Value = 1392.7407
$ws.Range("L31").Value = 4077.2727
$ws.Range("M31").Value = -1097.7407
$ws.Range("N31").Value = -4667.2727

$ws.Range("H34").Value = 2598.0408
$ws.Range("I34").Value = 1392.7407
$ws.Range("J34").Value = 4077.2727
$ws.Range("K34").Value = 1392.7407
$ws.Range("L34").Value = 4077.2727
$ws.Range("M34").Value = -1190.7407
$ws.Range("N34").Value = -4481.2727

$ws.Range("H58").Value = 4199.231
$ws.Range("I58").Value = 799.1667
$ws.Range("J58").Value = 45000
$ws.Range("K58").Value = 799.1667
$ws.Range("L58").Value = 45000
$ws.Range("M58").Value = -596.1667
$ws.Range("N58").Value = -45406

$ws.Range("H103").Value = 11500
$ws.Range("J103").Value = 15000
$ws.Range("L103").Value = 15000
$ws.Range("N103").Value = -17344

$ws.Range("H113").Value = 2981.6191
$ws.Range("I113").Value = 3293.4614
$ws.Range("J113").Value = 2474.875
$ws.Range("K113").Value = 3293.4614
$ws.Range("L113").Value = 2474.875
$ws.Range("M113").Value = -1123.4614
$ws.Range("N113").Value = -6814.875

$ws.Range("H136").Value = 4199.231
$ws.Range("I136").Value = 799.1667
$ws.Range("J136").Value = 45000
$ws.Range("K136").Value = 2397.5001
$ws.Range("L136").Value = 135000
$ws.Range("M136").Value = 152.4998999999998
$ws.Range("N136").Value = -140100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2097.7693
$ws.Range("I75").Value = 735.2
$ws.Range("J75").Value = 2949.375
$ws.Range("K75").Value = 2205.6
$ws.Range("L75").Value = 8848.125
$ws.Range("M75").Value = -1207.6
$ws.Range("N75").Value = -10844.125

$ws.Range("H78").Value = 2097.7693
$ws.Range("I78").Value = 735.2
$ws.Range("J78").Value = 2949.375
$ws.Range("K78").Value = 6616.8
$ws.Range("L78").Value = 26544.375
$ws.Range("M78").Value = -1624.8
$ws.Range("N78").Value = -36528.375

$ws.Range("H114").Value = 1300.4615
$ws.Range("I114").Value = 1881.6428
$ws.Range("J114").Value = 622.4167
$ws.Range("K114").Value = 5644.928400000001
$ws.Range("L114").Value = 1867.2501
$ws.Range("M114").Value = -2390.928400000001
$ws.Range("N114").Value = -8375.250099999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H80").Value = 4411.5557
$ws.Range("I80").Value = 7026.25
$ws.Range("J80").Value = 2319.8
$ws.Range("K80").Value = 7026.25
$ws.Range("L80").Value = 2319.8
$ws.Range("M80").Value = -6028.25
$ws.Range("N80").Value = -4315.8

$ws.Range("H83").Value = 4411.5557
$ws.Range("I83").Value = 7026.25
$ws.Range("J83").Value = 2319.8
$ws.Range("K83").Value = 35131.25
$ws.Range("L83").Value = 11599
$ws.Range("M83").Value = -30139.25
$ws.Range("N83").Value = -21583

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1667501
$ws.Range("I22").Value = 3030910
$ws.Range("J22").Value = 1112.1111
$ws.Range("K22").Value = 3030910
$ws.Range("L22").Value = 1112.1111
$ws.Range("M22").Value = -3030615
$ws.Range("N22").Value = -1702.1111

$ws.Range("H27").Value = 1667501
$ws.Range("I27").Value = 3030910
$ws.Range("J27").Value = 1112.1111
$ws.Range("K27").Value = 3030910
$ws.Range("L27").Value = 1112.1111
$ws.Range("M27").Value = -3030803
$ws.Range("N27").Value = -1326.1111

$ws.Range("H46").Value = 1518.1818
$ws.Range("I46").Value = 2200
$ws.Range("J46").Value = 700
$ws.Range("K46").Value = 2200
$ws.Range("L46").Value = 700
$ws.Range("M46").Value = -2012
$ws.Range("N46").Value = -1076

$ws.Range("H61").Value = 2719.8
$ws.Range("I61").Value = 2539.3333
$ws.Range("J61").Value = 2886.3845
$ws.Range("K61").Value = 2539.3333
$ws.Range("L61").Value = 2886.3845
$ws.Range("M61").Value = -2337.3333
$ws.Range("N61").Value = -3290.3845

$ws.Range("H113").Value = 2719.8
$ws.Range("I113").Value = 2539.3333
$ws.Range("J113").Value = 2886.3845
$ws.Range("K113").Value = 2539.3333
$ws.Range("L113").Value = 2886.3845
$ws.Range("M113").Value = -369.3332999999998
$ws.Range("N113").Value = -7226.3845

$ws.Range("H132").Value = 10876191
$ws.Range("I132").Value = 21749754
$ws.Range("J132").Value = 2626.4783
$ws.Range("K132").Value = 65249262
$ws.Range("L132").Value = 7879.4349
$ws.Range("M132").Value = -65246732
$ws.Range("N132").Value = -12939.4349

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2307.6428
$ws.Range("I132").Value = 1251
$ws.Range("J132").Value = 2730.3
$ws.Range("K132").Value = 3753
$ws.Range("L132").Value = 8190.900000000001
$ws.Range("M132").Value = -1223
$ws.Range("N132").Value = -13250.9

$ws.Range("H136").Value = 4337.0884
$ws.Range("I136").Value = 1052.1786
$ws.Range("J136").Value = 19666.666
$ws.Range("K136").Value = 3156.5358
$ws.Range("L136").Value = 58999.99800000001
$ws.Range("M136").Value = -606.5357999999997
$ws.Range("N136").Value = -64099.99800000001
